# Auto-generated edit script applying the Garuda_Profits.xlsx diff
# Updates cached numeric values for columns H-N across several leve rows
# on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1580
$ws.Cells.Item(40, 10).Value = 1836.3636
$ws.Cells.Item(40, 12).Value = 1836.3636
$ws.Cells.Item(40, 14).Value = -2186.3636
$ws.Cells.Item(64, 8).Value = 2947.946
$ws.Cells.Item(64, 9).Value = 2740
$ws.Cells.Item(64, 10).Value = 3024.963
$ws.Cells.Item(64, 11).Value = 2740
$ws.Cells.Item(64, 12).Value = 3024.963
$ws.Cells.Item(64, 13).Value = -2492
$ws.Cells.Item(64, 14).Value = -3520.963
$ws.Cells.Item(67, 8).Value = 2947.946
$ws.Cells.Item(67, 9).Value = 2740
$ws.Cells.Item(67, 10).Value = 3024.963
$ws.Cells.Item(67, 11).Value = 2740
$ws.Cells.Item(67, 12).Value = 3024.963
$ws.Cells.Item(67, 13).Value = -1882
$ws.Cells.Item(67, 14).Value = -4740.963
$ws.Cells.Item(74, 8).Value = 4314.615
$ws.Cells.Item(74, 9).Value = 4153.636
$ws.Cells.Item(74, 10).Value = 5200
$ws.Cells.Item(74, 11).Value = 4153.636
$ws.Cells.Item(74, 12).Value = 5200
$ws.Cells.Item(74, 13).Value = -3217.636
$ws.Cells.Item(74, 14).Value = -7072
$ws.Cells.Item(77, 8).Value = 4314.615
$ws.Cells.Item(77, 9).Value = 4153.636
$ws.Cells.Item(77, 10).Value = 5200
$ws.Cells.Item(77, 11).Value = 20768.18
$ws.Cells.Item(77, 12).Value = 26000
$ws.Cells.Item(77, 13).Value = -16088.18
$ws.Cells.Item(77, 14).Value = -35360
$ws.Cells.Item(106, 8).Value = 62751576
$ws.Cells.Item(106, 9).Value = 287516.44
$ws.Cells.Item(106, 10).Value = 500000000
$ws.Cells.Item(106, 11).Value = 287516.44
$ws.Cells.Item(106, 12).Value = 500000000
$ws.Cells.Item(106, 13).Value = -286885.44
$ws.Cells.Item(106, 14).Value = -500001262
$ws.Cells.Item(129, 8).Value = 523324.8
$ws.Cells.Item(129, 9).Value = 692.3333
$ws.Cells.Item(129, 10).Value = 549028.0600000001
$ws.Cells.Item(129, 11).Value = 2076.9999
$ws.Cells.Item(129, 12).Value = 1647084.18
$ws.Cells.Item(129, 13).Value = 2923.0001
$ws.Cells.Item(129, 14).Value = -1657084.18
$ws.Cells.Item(132, 8).Value = 3864240.2
$ws.Cells.Item(132, 9).Value = 4204790.5
$ws.Cells.Item(132, 10).Value = 4668.6665
$ws.Cells.Item(132, 11).Value = 12614371.5
$ws.Cells.Item(132, 12).Value = 14005.9995
$ws.Cells.Item(132, 13).Value = -12611841.5
$ws.Cells.Item(132, 14).Value = -19065.9995
$ws.Cells.Item(135, 8).Value = 804.5
$ws.Cells.Item(135, 9).Value = 533.2857
$ws.Cells.Item(135, 11).Value = 4799.571300000001
$ws.Cells.Item(135, 13).Value = -2264.571300000001
$ws.Cells.Item(137, 8).Value = 5883537
$ws.Cells.Item(137, 9).Value = 1245.6
$ws.Cells.Item(137, 11).Value = 3736.8
$ws.Cells.Item(137, 13).Value = -1186.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14921.964
$ws.Cells.Item(32, 9).Value = 15277.506
$ws.Cells.Item(32, 11).Value = 15277.506
$ws.Cells.Item(32, 13).Value = -14990.506
$ws.Cells.Item(38, 8).Value = 4000
$ws.Cells.Item(38, 9).Value = 4000
$ws.Cells.Item(38, 11).Value = 4000
$ws.Cells.Item(38, 13).Value = -3533
$ws.Cells.Item(74, 8).Value = 585.03845
$ws.Cells.Item(74, 9).Value = 441.4091
$ws.Cells.Item(74, 11).Value = 441.4091
$ws.Cells.Item(74, 13).Value = 432.5909
$ws.Cells.Item(77, 8).Value = 585.03845
$ws.Cells.Item(77, 9).Value = 441.4091
$ws.Cells.Item(77, 11).Value = 2207.0455
$ws.Cells.Item(77, 13).Value = 2160.9545
$ws.Cells.Item(97, 8).Value = 942.8333
$ws.Cells.Item(97, 9).Value = 942.8333
$ws.Cells.Item(97, 11).Value = 942.8333
$ws.Cells.Item(97, 13).Value = -446.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1241.8572
$ws.Cells.Item(86, 9).Value = 1138.6
$ws.Cells.Item(86, 10).Value = 1500
$ws.Cells.Item(86, 11).Value = 1138.6
$ws.Cells.Item(86, 12).Value = 1500
$ws.Cells.Item(86, 13).Value = -15.59999999999991
$ws.Cells.Item(86, 14).Value = -3746
$ws.Cells.Item(89, 8).Value = 1241.8572
$ws.Cells.Item(89, 9).Value = 1138.6
$ws.Cells.Item(89, 10).Value = 1500
$ws.Cells.Item(89, 11).Value = 5693
$ws.Cells.Item(89, 12).Value = 7500
$ws.Cells.Item(89, 13).Value = -77
$ws.Cells.Item(89, 14).Value = -18732
$ws.Cells.Item(93, 8).Value = 30000
$ws.Cells.Item(93, 10).Value = 30000
$ws.Cells.Item(93, 12).Value = 30000
$ws.Cells.Item(93, 14).Value = -33744
$ws.Cells.Item(134, 8).Value = 30237.611
$ws.Cells.Item(134, 9).Value = 46067.434
$ws.Cells.Item(134, 10).Value = 2231
$ws.Cells.Item(134, 11).Value = 138202.302
$ws.Cells.Item(134, 12).Value = 6693
$ws.Cells.Item(134, 13).Value = -135667.302
$ws.Cells.Item(134, 14).Value = -11763

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6252560.5
$ws.Cells.Item(31, 9).Value = 2807.5264
$ws.Cells.Item(31, 10).Value = 15386815
$ws.Cells.Item(31, 11).Value = 2807.5264
$ws.Cells.Item(31, 12).Value = 15386815
$ws.Cells.Item(31, 13).Value = -2512.5264
$ws.Cells.Item(31, 14).Value = -15387405
$ws.Cells.Item(34, 8).Value = 6252560.5
$ws.Cells.Item(34, 9).Value = 2807.5264
$ws.Cells.Item(34, 10).Value = 15386815
$ws.Cells.Item(34, 11).Value = 2807.5264
$ws.Cells.Item(34, 12).Value = 15386815
$ws.Cells.Item(34, 13).Value = -2605.5264
$ws.Cells.Item(34, 14).Value = -15387219
$ws.Cells.Item(132, 8).Value = 2464.0625
$ws.Cells.Item(132, 9).Value = 2259.6191
$ws.Cells.Item(132, 10).Value = 2854.3635
$ws.Cells.Item(132, 11).Value = 6778.8573
$ws.Cells.Item(132, 12).Value = 8563.0905
$ws.Cells.Item(132, 13).Value = -4248.8573
$ws.Cells.Item(132, 14).Value = -13623.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1350.5647
$ws.Cells.Item(68, 9).Value = 1300.6666
$ws.Cells.Item(68, 10).Value = 1392.8695
$ws.Cells.Item(68, 11).Value = 3901.9998
$ws.Cells.Item(68, 12).Value = 4178.6085
$ws.Cells.Item(68, 13).Value = -3090.9998
$ws.Cells.Item(68, 14).Value = -5800.6085
$ws.Cells.Item(71, 8).Value = 1350.5647
$ws.Cells.Item(71, 9).Value = 1300.6666
$ws.Cells.Item(71, 10).Value = 1392.8695
$ws.Cells.Item(71, 11).Value = 11705.9994
$ws.Cells.Item(71, 12).Value = 12535.8255
$ws.Cells.Item(71, 13).Value = -7649.999400000001
$ws.Cells.Item(71, 14).Value = -20647.8255
$ws.Cells.Item(75, 8).Value = 956.5
$ws.Cells.Item(75, 9).Value = 956.5
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 2869.5
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = -1871.5
$ws.Cells.Item(75, 14).Value = $null
$ws.Cells.Item(78, 8).Value = 956.5
$ws.Cells.Item(78, 9).Value = 956.5
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 8608.5
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 13).Value = -3616.5
$ws.Cells.Item(78, 14).Value = $null
$ws.Cells.Item(113, 8).Value = 645.0192
$ws.Cells.Item(113, 9).Value = 508.4074
$ws.Cells.Item(113, 10).Value = 792.5599999999999
$ws.Cells.Item(113, 11).Value = 1525.2222
$ws.Cells.Item(113, 12).Value = 2377.68
$ws.Cells.Item(113, 13).Value = 644.7778000000001
$ws.Cells.Item(113, 14).Value = -6717.68
$ws.Cells.Item(122, 8).Value = 743.7059
$ws.Cells.Item(122, 9).Value = 838.12
$ws.Cells.Item(122, 10).Value = 652.9231
$ws.Cells.Item(122, 11).Value = 7543.08
$ws.Cells.Item(122, 12).Value = 5876.3079
$ws.Cells.Item(122, 13).Value = -5093.08
$ws.Cells.Item(122, 14).Value = -10776.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2390.0908
$ws.Cells.Item(126, 9).Value = 1049.2
$ws.Cells.Item(126, 10).Value = 3507.5
$ws.Cells.Item(126, 11).Value = 3147.6
$ws.Cells.Item(126, 12).Value = 10522.5
$ws.Cells.Item(126, 13).Value = -677.6000000000004
$ws.Cells.Item(126, 14).Value = -15462.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 581.2857
$ws.Cells.Item(16, 9).Value = 581.2857
$ws.Cells.Item(16, 11).Value = 581.2857
$ws.Cells.Item(16, 13).Value = -411.2857
$ws.Cells.Item(40, 8).Value = 127226
$ws.Cells.Item(40, 9).Value = 144829.72
$ws.Cells.Item(40, 10).Value = 4000
$ws.Cells.Item(40, 11).Value = 144829.72
$ws.Cells.Item(40, 12).Value = 4000
$ws.Cells.Item(40, 13).Value = -144693.72
$ws.Cells.Item(40, 14).Value = -4272
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).Value = $null
$ws.Cells.Item(61, 8).Value = 4127
$ws.Cells.Item(61, 9).Value = 3304
$ws.Cells.Item(61, 10).Value = 4950
$ws.Cells.Item(61, 11).Value = 3304
$ws.Cells.Item(61, 12).Value = 4950
$ws.Cells.Item(61, 13).Value = -3102
$ws.Cells.Item(61, 14).Value = -5354
$ws.Cells.Item(93, 8).Value = 1127284.1
$ws.Cells.Item(93, 9).Value = 1502556.6
$ws.Cells.Item(93, 11).Value = 1502556.6
$ws.Cells.Item(93, 13).Value = -1501308.6
$ws.Cells.Item(100, 8).Value = 4331.3335
$ws.Cells.Item(100, 9).Value = 4000
$ws.Cells.Item(100, 11).Value = 4000
$ws.Cells.Item(100, 13).Value = -3459
$ws.Cells.Item(113, 8).Value = 4127
$ws.Cells.Item(113, 9).Value = 3304
$ws.Cells.Item(113, 10).Value = 4950
$ws.Cells.Item(113, 11).Value = 3304
$ws.Cells.Item(113, 12).Value = 4950
$ws.Cells.Item(113, 13).Value = -1134
$ws.Cells.Item(113, 14).Value = -9290
$ws.Cells.Item(132, 8).Value = 13621.167
$ws.Cells.Item(132, 9).Value = 20234.092
$ws.Cells.Item(132, 10).Value = 3229.4285
$ws.Cells.Item(132, 11).Value = 60702.276
$ws.Cells.Item(132, 12).Value = 9688.2855
$ws.Cells.Item(132, 13).Value = -58172.276
$ws.Cells.Item(132, 14).Value = -14748.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1485.7142
$ws.Cells.Item(81, 9).Value = 1280
$ws.Cells.Item(81, 10).Value = 2000
$ws.Cells.Item(81, 11).Value = 2560
$ws.Cells.Item(81, 12).Value = 4000
$ws.Cells.Item(81, 13).Value = -1499
$ws.Cells.Item(81, 14).Value = -6122
$ws.Cells.Item(84, 8).Value = 1485.7142
$ws.Cells.Item(84, 9).Value = 1280
$ws.Cells.Item(84, 10).Value = 2000
$ws.Cells.Item(84, 11).Value = 12800
$ws.Cells.Item(84, 12).Value = 20000
$ws.Cells.Item(84, 13).Value = -7496
$ws.Cells.Item(84, 14).Value = -30608

Write-Output "Applied 234 cell updates"